# Auto-generated edit script applying the Typhon_Profits diff
# Updates per-leve currentAveragePrice / LevePrice / LeveProfit figures across 8 job sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 168.42857
$ws.Range("I4").Value = 168.42857
$ws.Range("K4").Value = 168.42857
$ws.Range("M4").Value = -54.42857000000001
$ws.Range("H55").Value = 933.3333
$ws.Range("I55").Value = 1733.3334
$ws.Range("J55").Value = 133.33333
$ws.Range("K55").Value = 1733.3334
$ws.Range("L55").Value = 133.33333
$ws.Range("M55").Value = -1519.3334
$ws.Range("N55").Value = -561.3333299999999
$ws.Range("H138").Value = 2676.7715
$ws.Range("I138").Value = 1658.2354
$ws.Range("J138").Value = 3638.7222
$ws.Range("K138").Value = 4974.706200000001
$ws.Range("L138").Value = 10916.1666
$ws.Range("M138").Value = 165.2937999999995
$ws.Range("N138").Value = -21196.1666
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 340
$ws.Range("I4").Value = 180
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 180
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -64
$ws.Range("N4").Value = -732
$ws.Range("H61").Value = 3446.2727
$ws.Range("I61").Value = 3217.389
$ws.Range("J61").Value = 3720.9333
$ws.Range("K61").Value = 3217.389
$ws.Range("L61").Value = 3720.9333
$ws.Range("M61").Value = -3005.389
$ws.Range("N61").Value = -4144.933300000001
$ws.Range("H102").Value = 7800
$ws.Range("I102").Value = 6333.3335
$ws.Range("K102").Value = 6333.3335
$ws.Range("M102").Value = -4711.3335
$ws.Range("H136").Value = 3446.2727
$ws.Range("I136").Value = 3217.389
$ws.Range("J136").Value = 3720.9333
$ws.Range("K136").Value = 9652.167000000001
$ws.Range("L136").Value = 11162.7999
$ws.Range("M136").Value = -7102.167000000001
$ws.Range("N136").Value = -16262.7999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 8214
$ws.Range("I97").Value = 4656.8
$ws.Range("J97").Value = 26000
$ws.Range("K97").Value = 4656.8
$ws.Range("L97").Value = 26000
$ws.Range("M97").Value = -3665.8
$ws.Range("N97").Value = -27982
$ws.Range("H99").Value = 1333.3334
$ws.Range("I99").Value = 1333.3334
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1333.3334
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 164.6666
$ws.Range("N99").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 46.416668
$ws.Range("I7").Value = 46.090908
$ws.Range("K7").Value = 46.090908
$ws.Range("M7").Value = 66.909092
$ws.Range("H16").Value = 2800
$ws.Range("I16").Value = 2800
$ws.Range("K16").Value = 2800
$ws.Range("M16").Value = -2513
$ws.Range("H31").Value = 14477.148
$ws.Range("I31").Value = 31344
$ws.Range("K31").Value = 31344
$ws.Range("M31").Value = -31049
$ws.Range("H34").Value = 14477.148
$ws.Range("I34").Value = 31344
$ws.Range("K34").Value = 31344
$ws.Range("M34").Value = -31142
$ws.Range("H113").Value = 2800
$ws.Range("I113").Value = 2800
$ws.Range("K113").Value = 2800
$ws.Range("M113").Value = -630
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 841.6539
$ws.Range("I5").Value = 950.375
$ws.Range("K5").Value = 2851.125
$ws.Range("M5").Value = -2739.125
$ws.Range("H33").Value = 299.66666
$ws.Range("I33").Value = 299
$ws.Range("K33").Value = 1794
$ws.Range("M33").Value = -1511
$ws.Range("H68").Value = 4182.4116
$ws.Range("J68").Value = 6151.5454
$ws.Range("L68").Value = 18454.6362
$ws.Range("N68").Value = -20076.6362
$ws.Range("H71").Value = 4182.4116
$ws.Range("J71").Value = 6151.5454
$ws.Range("L71").Value = 55363.9086
$ws.Range("N71").Value = -63475.9086
$ws.Range("H107").Value = 3183.077
$ws.Range("J107").Value = 667.6429000000001
$ws.Range("L107").Value = 2002.9287
$ws.Range("N107").Value = -5842.9287
$ws.Range("H131").Value = 791.14
$ws.Range("J131").Value = 808.28864
$ws.Range("L131").Value = 2424.86592
$ws.Range("N131").Value = -12504.86592
$ws.Range("H132").Value = 897.8
$ws.Range("J132").Value = 500
$ws.Range("L132").Value = 4500
$ws.Range("N132").Value = -9560
$ws.Range("H135").Value = 841.6539
$ws.Range("I135").Value = 950.375
$ws.Range("K135").Value = 8553.375
$ws.Range("M135").Value = -6018.375
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5547.2812
$ws.Range("I126").Value = 5317.5
$ws.Range("J126").Value = 5930.25
$ws.Range("K126").Value = 15952.5
$ws.Range("L126").Value = 17790.75
$ws.Range("M126").Value = -13482.5
$ws.Range("N126").Value = -22730.75
$ws.Range("H132").Value = 18515.848
$ws.Range("J132").Value = 52839.8
$ws.Range("L132").Value = 158519.4
$ws.Range("N132").Value = -163579.4
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3600
$ws.Range("I22").Value = 3600
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3600
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -3305
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 3600
$ws.Range("I27").Value = 3600
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 3600
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -3493
$ws.Range("N27").ClearContents()
$ws.Range("H40").Value = 373670.34
$ws.Range("I40").Value = 373670.34
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 373670.34
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -373534.34
$ws.Range("N40").ClearContents()
$ws.Range("H46").Value = 1727.2354
$ws.Range("I46").Value = 1804.5
$ws.Range("J46").Value = 1366.6666
$ws.Range("K46").Value = 1804.5
$ws.Range("L46").Value = 1366.6666
$ws.Range("M46").Value = -1616.5
$ws.Range("N46").Value = -1742.6666
$ws.Range("H132").Value = 1670.9117
$ws.Range("I132").Value = 1148.4762
$ws.Range("J132").Value = 2514.8462
$ws.Range("K132").Value = 3445.4286
$ws.Range("L132").Value = 7544.5386
$ws.Range("M132").Value = -915.4286000000002
$ws.Range("N132").Value = -12604.5386
$ws.Range("H136").Value = 14966.553
$ws.Range("I136").Value = 21756.666
$ws.Range("J136").Value = 3326.3572
$ws.Range("K136").Value = 65269.99800000001
$ws.Range("L136").Value = 9979.071599999999
$ws.Range("M136").Value = -62719.99800000001
$ws.Range("N136").Value = -15079.0716
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2080099.5
$ws.Range("I113").Value = 1474.3334
$ws.Range("J113").Value = 6757006
$ws.Range("K113").Value = 4423.0002
$ws.Range("L113").Value = 20271018
$ws.Range("M113").Value = -2253.0002
$ws.Range("N113").Value = -20275358
$ws.Range("H122").Value = 1865.9474
$ws.Range("I122").Value = 1691.625
$ws.Range("K122").Value = 5074.875
$ws.Range("M122").Value = -2624.875
$ws.Range("H126").Value = 1524.75
$ws.Range("I126").Value = 649.5
$ws.Range("J126").Value = 2400
$ws.Range("K126").Value = 1948.5
$ws.Range("L126").Value = 7200
$ws.Range("M126").Value = 521.5
$ws.Range("N126").Value = -12140
$ws.Range("H132").Value = 2118.0833
$ws.Range("I132").Value = 1826.8
$ws.Range("J132").Value = 3574.5
$ws.Range("K132").Value = 5480.4
$ws.Range("L132").Value = 10723.5
$ws.Range("M132").Value = -2950.4
$ws.Range("N132").Value = -15783.5
$ws.Range("H136").Value = 1258.7241
$ws.Range("I136").Value = 825.5
$ws.Range("K136").Value = 2476.5
$ws.Range("M136").Value = 73.5

Write-Host "Applied Typhon_Profits updates"
